# Applies the "Updated symbol list" edit (Fri Dec 23 02:55:01 UTC 2022).
# Cryptocurrency rows were reshuffled/updated and several Price cells were
# refreshed; every value in this sheet is stored as literal text (t="inlineStr"
# in the source file), so numeric-looking prices are written with a leading
# apostrophe to force Excel to keep them as text instead of coercing them to
# real numbers (which would silently drop significant trailing/leading zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = [ordered]@{
    'D2' = '''246.13'
    'D3' = '''21.98'
    'D4' = '''5.419'
    'D5' = '''0.05764'
    'D7' = '''6.328'
    'D8' = '''0.8133'
    'D9' = '''0.9852'
    'E9' = '8FTXTokenFTT'
    'B10' = 'One'
    'C10' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D10' = '''0.01123'
    'E10' = '9OneONEBestin24h'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '''0.1432'
    'E11' = '10WazirXWRX'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '''0.07412'
    'E12' = '11MandalaExchangeTokenMDX'
    'B13' = 'LiechtensteinCryptoassetsExchange'
    'C13' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D13' = '''0.03146'
    'E13' = '12LiechtensteinCryptoassetsExchangeLCX'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D14' = '''0.02996'
    'E14' = '13BitrueCoinBTR'
    'B15' = 'MCDex'
    'C15' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D15' = '''4.148'
    'E15' = '14MCDexMCB'
    'B16' = 'BitMartToken'
    'C16' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D16' = '''0.09403'
    'E16' = '15BitMartTokenBMX'
    'B17' = 'BitForexToken'
    'C17' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D17' = '''0.001595'
    'E17' = '16BitForexTokenBF'
    'B18' = 'CoinExToken'
    'C18' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D18' = '''0.04811'
    'E18' = '17CoinExTokenCET'
    'D19' = '''0.006178'
    'D20' = '''0.004119'
    'D21' = '''0.0009953'
    'D22' = '''0.0001499'
    'D23' = '''3.756'
    'D24' = '''2.208'
    'D25' = '''0.3258'
    'D26' = '''0.1330'
    'D27' = '''0.0003996'
    'D40' = '''0.03888'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D41' = '''0.006456'
    'E41' = '40KickTokenKICK'
    'B42' = 'BKEXToken'
    'C42' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D42' = '''0.1072'
    'E42' = '41BKEXTokenBKK'
    'B43' = 'CEJI'
    'C43' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D43' = '''0.002998'
    'E43' = '42CEJICEJI'
    'D44' = '''0.006271'
    'D45' = '''0.00005589'
    'D46' = '''0.00000000749'
    'D47' = '''0.3797'
    'D50' = '''0.01009'
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}
